$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Admin")

# Row 9: B9 now shows the PA user's email (new value) and gets its own
# mailto hyperlink. C9 keeps showing "Test@123" (unchanged, already wired
# to its existing hyperlink rId3).
$ws.Range("B9").Value = "pauser@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:pauser@gmail.com")

# Row 10: a new data row mirroring row 9 - "Valid User PA Name" label in A,
# the old PA email in B, and "Test@123" (the password) in C - each of B10
# and C10 getting their own mailto hyperlink.
$ws.Range("A10").Value = "Valid User PA Name"

$ws.Range("B10").Value = "testuserpa@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:testuserpa@gmail.com")

$ws.Range("C10").Value = "Test@123"
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:Test@123")

# Reflect the author's last selection before saving.
[void]$ws.Activate()
[void]$ws.Range("C14").Select()
